# Auto-update gym prices

$wb = $excel.ActiveWorkbook

# Helper: set a cell's value as literal text, avoiding Excel's automatic
# number/currency conversion (e.g. turning "$2,048.00" into a numeric
# currency value). We briefly switch the cell to Text format, assign the
# value, then restore the original number format so the cell's appearance
# is unchanged.
function Set-LiteralText {
    param($range, [string]$text)
    $origFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $origFormat
}

# Sheet 1: "4x4 Squat Racks"
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
Set-LiteralText $ws1.Range("C2") "$2,048.00"
Set-LiteralText $ws1.Range("C3") "Not available"

# Sheet 2: "Squat Stands"
$ws2 = $wb.Worksheets.Item("Squat Stands")
Set-LiteralText $ws2.Range("C2") "$1,481.00"
Set-LiteralText $ws2.Range("C3") "Price not available"
